$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.423.69'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.563.64'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '287.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3651'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.76%  '
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("E9").Value = '  -1.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.127'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07408'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.15%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.93'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.937'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.875'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.23%  '
$ws.Range("D16").Value = '1.563.40'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("E17").Value = '  -1.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06735'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.297'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.26%  '
$ws.Range("D24").Value = '22.407.90'
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.371'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.551'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.05%  '
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.32%  '
$ws.Range("D31").Value = '1.739.59'
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.051'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.098'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.988'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.622'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08271'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02392'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.311'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2217'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06376'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.322'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6089'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.44%  '
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.768'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5736'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("E48").Value = '  -4.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.227'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07244'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.19%  '
